$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = '2025-12-17 Wednesday'

# Update each table cell value (positional replacement; some old values repeat)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = '76-31='
$t.Cell(1,2).Range.Text = '52-43='
$t.Cell(1,3).Range.Text = '49-39='
$t.Cell(1,4).Range.Text = '32+10='
$t.Cell(1,5).Range.Text = '32-5='
$t.Cell(2,1).Range.Text = '74-45='
$t.Cell(2,2).Range.Text = '33-24='
$t.Cell(2,3).Range.Text = '69-43='
$t.Cell(2,4).Range.Text = '45-4='
$t.Cell(2,5).Range.Text = '57-47='
$t.Cell(3,1).Range.Text = '84-44='
$t.Cell(3,2).Range.Text = '67-2='
$t.Cell(3,3).Range.Text = '2+4='
$t.Cell(3,4).Range.Text = '79-16='
$t.Cell(3,5).Range.Text = '21-8='
$t.Cell(4,1).Range.Text = '8+68='
$t.Cell(4,2).Range.Text = '17+0='
$t.Cell(4,3).Range.Text = '58+13='
$t.Cell(4,4).Range.Text = '38+9='
$t.Cell(4,5).Range.Text = '85-43='
$t.Cell(5,1).Range.Text = '13-9='
$t.Cell(5,2).Range.Text = '14+58='
$t.Cell(5,3).Range.Text = '78-10='
$t.Cell(5,4).Range.Text = '83-68='
$t.Cell(5,5).Range.Text = '47+9='
$t.Cell(6,1).Range.Text = '26+50='
$t.Cell(6,2).Range.Text = '6+37='
$t.Cell(6,3).Range.Text = '86+13='
$t.Cell(6,4).Range.Text = '11+35='
$t.Cell(6,5).Range.Text = '23+70='
$t.Cell(7,1).Range.Text = '64+23='
$t.Cell(7,2).Range.Text = '17+82='
$t.Cell(7,3).Range.Text = '41+42='
$t.Cell(7,4).Range.Text = '88-27='
$t.Cell(7,5).Range.Text = '52+35='
$t.Cell(8,1).Range.Text = '89-15='
$t.Cell(8,2).Range.Text = '4+94='
$t.Cell(8,3).Range.Text = '9+19='
$t.Cell(8,4).Range.Text = '81-71='
$t.Cell(8,5).Range.Text = '94-57='
$t.Cell(9,1).Range.Text = '91-20='
$t.Cell(9,2).Range.Text = '98-35='
$t.Cell(9,3).Range.Text = '88-18='
$t.Cell(9,4).Range.Text = '46-44='
$t.Cell(9,5).Range.Text = '63-24='
$t.Cell(10,1).Range.Text = '51+19='
$t.Cell(10,2).Range.Text = '41+58='
$t.Cell(10,3).Range.Text = '48-41='
$t.Cell(10,4).Range.Text = '25-10='
$t.Cell(10,5).Range.Text = '33+56='
$t.Cell(11,1).Range.Text = '11+5='
$t.Cell(11,2).Range.Text = '96-62='
$t.Cell(11,3).Range.Text = '13+73='
$t.Cell(11,4).Range.Text = '71-59='
$t.Cell(11,5).Range.Text = '97-77='
$t.Cell(12,1).Range.Text = '70-63='
$t.Cell(12,2).Range.Text = '82-82='
$t.Cell(12,3).Range.Text = '22-18='
$t.Cell(12,4).Range.Text = '11+36='
$t.Cell(12,5).Range.Text = '48+49='
$t.Cell(13,1).Range.Text = '82-18='
$t.Cell(13,2).Range.Text = '25-16='
$t.Cell(13,3).Range.Text = '23+60='
$t.Cell(13,4).Range.Text = '38+54='
$t.Cell(13,5).Range.Text = '8+68='
$t.Cell(14,1).Range.Text = '24+1='
$t.Cell(14,2).Range.Text = '83-25='
$t.Cell(14,3).Range.Text = '58-8='
$t.Cell(14,4).Range.Text = '35-15='
$t.Cell(14,5).Range.Text = '80+10='
$t.Cell(15,1).Range.Text = '45-30='
$t.Cell(15,2).Range.Text = '24+59='
$t.Cell(15,3).Range.Text = '46-13='
$t.Cell(15,4).Range.Text = '57-25='
$t.Cell(15,5).Range.Text = '43+22='
$t.Cell(16,1).Range.Text = '54+37='
$t.Cell(16,2).Range.Text = '32+55='
$t.Cell(16,3).Range.Text = '50-0='
$t.Cell(16,4).Range.Text = '15+73='
$t.Cell(16,5).Range.Text = '83-43='
$t.Cell(17,1).Range.Text = '8+65='
$t.Cell(17,2).Range.Text = '55-5='
$t.Cell(17,3).Range.Text = '51+12='
$t.Cell(17,4).Range.Text = '67+6='
$t.Cell(17,5).Range.Text = '52-30='
$t.Cell(18,1).Range.Text = '86-30='
$t.Cell(18,2).Range.Text = '66+30='
$t.Cell(18,3).Range.Text = '9+55='
$t.Cell(18,4).Range.Text = '51+33='
$t.Cell(18,5).Range.Text = '72-28='
$t.Cell(19,1).Range.Text = '87-85='
$t.Cell(19,2).Range.Text = '52+28='
$t.Cell(19,3).Range.Text = '74-47='
$t.Cell(19,4).Range.Text = '94-91='
$t.Cell(19,5).Range.Text = '32+23='
$t.Cell(20,1).Range.Text = '63-9='
$t.Cell(20,2).Range.Text = '47+47='
$t.Cell(20,3).Range.Text = '66+25='
$t.Cell(20,4).Range.Text = '33+53='
$t.Cell(20,5).Range.Text = '63+34='
